$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the simulation object name references (shared strings) used in the
# XY plot table: "PVSystem.bb_258064_1_5" -> "PVSystem.oh_263163_0_1"
# and "dev_258183_131" -> "dev_263265_0_1"
$ws.Range("B3").Value = "PVSystem.oh_263163_0_1"
$ws.Range("G3").Value = "PVSystem.oh_263163_0_1"
$ws.Range("B4").Value = "PVSystem.oh_263163_0_1"
$ws.Range("G4").Value = "dev_263265_0_1"

# Update the active cell selection on the sheet view
$ws.Range("F14").Select()
